$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each target cell holds plain text (inlineStr) in the workbook, including
# numeric-looking price strings that must keep exact trailing zeros (e.g. "0.06940").
# Force text entry via a temporary Text number format so Excel does not
# auto-convert the literal into a Number, then clear the format again so the
# cell style index is left exactly as it was before the edit.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "244.13"
$ws.Range("D2").ClearFormats()
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "24.56"
$ws.Range("D3").ClearFormats()
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "LEO"
$ws.Range("B4").ClearFormats()
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C4").ClearFormats()
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "3.501"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "3LEOLEO"
$ws.Range("E4").ClearFormats()
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "HuobiToken"
$ws.Range("B5").ClearFormats()
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("C5").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "5.119"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4HuobiTokenHT"
$ws.Range("E5").ClearFormats()
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "Cronos"
$ws.Range("B6").ClearFormats()
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("C6").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.05753"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "5CronosCRO"
$ws.Range("E6").ClearFormats()
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("B7").ClearFormats()
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("C7").ClearFormats()
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.470"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "6KuCoinTokenKCS"
$ws.Range("E7").ClearFormats()
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "GateToken"
$ws.Range("B8").ClearFormats()
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("C8").ClearFormats()
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.120"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "7GateTokenGT"
$ws.Range("E8").ClearFormats()
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "MXToken"
$ws.Range("B9").ClearFormats()
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C9").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8103"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "8MXTokenMX"
$ws.Range("E9").ClearFormats()
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "FTXToken"
$ws.Range("B10").ClearFormats()
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("C10").ClearFormats()
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8415"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "9FTXTokenFTT"
$ws.Range("E10").ClearFormats()
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "WazirX"
$ws.Range("B11").ClearFormats()
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("C11").ClearFormats()
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1337"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("E11").ClearFormats()
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("B12").ClearFormats()
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("C12").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06940"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("E12").ClearFormats()
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("B13").ClearFormats()
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("C13").ClearFormats()
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03134"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("E13").ClearFormats()
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("B14").ClearFormats()
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("C14").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.02844"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("E14").ClearFormats()
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("B15").ClearFormats()
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("C15").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09366"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("E15").ClearFormats()
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "MCDex"
$ws.Range("B16").ClearFormats()
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("C16").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.762"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("E16").ClearFormats()
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("B17").ClearFormats()
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("C17").ClearFormats()
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001513"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("E17").ClearFormats()
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("B18").ClearFormats()
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("C18").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04646"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("E18").ClearFormats()
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "One"
$ws.Range("B19").ClearFormats()
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("C19").ClearFormats()
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0006009"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "18OneONEWorstin24h"
$ws.Range("E19").ClearFormats()
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "TigerCash"
$ws.Range("B20").ClearFormats()
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("C20").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.006112"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "19TigerCashTCH"
$ws.Range("E20").ClearFormats()
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "BitKan"
$ws.Range("B21").ClearFormats()
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("C21").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001234"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "20BitKanKAN"
$ws.Range("E21").ClearFormats()
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "HotbitToken"
$ws.Range("B22").ClearFormats()
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("C22").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.004278"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "21HotbitTokenHTB"
$ws.Range("E22").ClearFormats()
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "NitroEx"
$ws.Range("B23").ClearFormats()
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("C23").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.00008695"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "22NitroExNTX"
$ws.Range("E23").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.082"
$ws.Range("D24").ClearFormats()
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1358"
$ws.Range("D27").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0002326"
$ws.Range("D28").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03613"
$ws.Range("D40").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006416"
$ws.Range("D41").ClearFormats()
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002930"
$ws.Range("D43").ClearFormats()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007382"
$ws.Range("D44").ClearFormats()
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005298"
$ws.Range("D45").ClearFormats()
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.2510"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
$ws.Range("E47").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002274"
$ws.Range("D48").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001999"
$ws.Range("D50").ClearFormats()
